$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B4: Cell_Count now references the named range Pack_Count instead of the raw cell B3
$ws.Range("B4").Formula = "=Pack_Count*3"

# B5: Battery_Capacity formula corrected from mAh(*5000) to Ah(*5) basis, also using Pack_Count
$ws.Range("B5").Formula = "=Pack_Count*5*11.1*3600"

# B6: Battery_Weight now references the named range Pack_Count instead of the raw cell B3
$ws.Range("B6").Formula = "=Pack_Count*3*69"

# E5: fix note text from mAh to Ah
$ws.Range("E5").Value = '"=capacity(Ah)*11.1(Vnominal)*36000(sec/hour)'

# H2: touch the alignment so a new cell style (applyAlignment) gets stamped on the notes cell
$ws.Range("H2").HorizontalAlignment = 1

# Column B: widen slightly to fit the updated values
$ws.Columns("B").ColumnWidth = 11.67

# Final selection left on E5 (last cell touched)
[void]$ws.Range("E5").Select()
